# Normalize the "North Ireland" sheet header row so it matches the
# "Wales" sheet's header row (labels + formatting), now that the
# workbook needs to support multiple sheets consistently.

$wb = $excel.ActiveWorkbook
$wales = $wb.Worksheets.Item("Wales")
$ni = $wb.Worksheets.Item("North Ireland")

# Update the North Ireland header labels to match the common labels
# used on the Wales sheet.
$ni.Range("A1").Value = "Location"
$ni.Range("D1").Value = "Likelihood (%)"
$ni.Range("E1").Value = "Recommended Reason"

# Copy the header row formatting from Wales onto North Ireland so both
# sheets share the same header style.
$wales.Range("A1:F1").Copy()
$ni.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reset the view on the Wales sheet: select the header row.
$wales.Select()
$wales.Range("A1:F1").Select()

# Reset the view on the North Ireland sheet: scroll back to the top and
# select the header row (previously it was scrolled down to A4 with
# just A4 selected).
$ni.Select()
$ni.Range("A1:F1").Select()

$wb.Windows.Item(1).Width = 12375
